# Iraq League workbook update
# - Re-orders several pairs of match rows that were out of chronological /
#   id order by swapping their full B:AC content (this also fixes a
#   "Al Sinaah" / "Al Karkh" team-name mix-up that only affected those
#   particular rows).
# - Appends new match results, turning the previous last row into the
#   3rd-to-last, and filling in its previously-missing result/odds data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# A handful of rows were recorded in the wrong order relative to
# their neighbour; swap each pair's full B:AC payload (id in column
# A stays put, everything else moves together).
# ---------------------------------------------------------------------
$pairs = @(
    @(4, 5),
    @(17, 18),
    @(22, 23),
    @(41, 42),
    @(69, 70),
    @(73, 74),
    @(78, 79),
    @(108, 109)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $rangeA = $ws.Range("B$r1`:AC$r1")
    $rangeB = $ws.Range("B$r2`:AC$r2")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# ---------------------------------------------------------------------
# Insert two fresh rows ahead of the final match (old row 147, which
# becomes row 149), then populate the three rows with the latest
# results.
# ---------------------------------------------------------------------
$ws.Cells.Item(147, 1).EntireRow.Insert()
$ws.Cells.Item(147, 1).EntireRow.Insert()

# Carry over the A/E column formatting (bold+border id style, date
# number format) from the row above onto the two new rows.
$ws.Range("A146").Copy()
$ws.Range("A147:A148").PasteSpecial(-4122)
$ws.Range("E146").Copy()
$ws.Range("E147:E148").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 147 (new match)
$ws.Range("A147").Value2 = 145
$ws.Range("B147").Value2 = 7935060
$ws.Range("C147").Value2 = "Iraq League"
$ws.Range("D147").Value2 = "Iraq League"
$ws.Range("E147").Value2 = 45360.35416666666
$ws.Range("F147").Value2 = "Al Quwa Al Jawiya"
$ws.Range("G147").Value2 = "Amanat Baghdad"
$ws.Range("H147").Value2 = 1
$ws.Range("I147").Value2 = 1
$ws.Range("J147").Value2 = "D"
$ws.Range("K147").Value2 = 1.333
$ws.Range("L147").Value2 = 4.333
$ws.Range("M147").Value2 = 7.5
$ws.Range("N147").Value2 = 1.285
$ws.Range("O147").Value2 = 4.333
$ws.Range("P147").Value2 = 10
$ws.Range("Q147").Value2 = -1.75
$ws.Range("R147").Value2 = 1.95
$ws.Range("S147").Value2 = 1.85
$ws.Range("T147").Value2 = 2.75
$ws.Range("U147").Value2 = 1.95
$ws.Range("V147").Value2 = 1.85
$ws.Range("W147").Value2 = -1
$ws.Range("X147").Value2 = 3.333
$ws.Range("Y147").Value2 = -1
$ws.Range("Z147").Value2 = -1
$ws.Range("AA147").Value2 = 0.8500000000000001
$ws.Range("AB147").Value2 = -1
$ws.Range("AC147").Value2 = 0.8500000000000001

# Row 148 (new match)
$ws.Range("A148").Value2 = 146
$ws.Range("B148").Value2 = 7935122
$ws.Range("C148").Value2 = "Iraq League"
$ws.Range("D148").Value2 = "Iraq League"
$ws.Range("E148").Value2 = 45360.45833333334
$ws.Range("F148").Value2 = "Zakho"
$ws.Range("G148").Value2 = "Naft AlWasat"
$ws.Range("H148").Value2 = 0
$ws.Range("I148").Value2 = 0
$ws.Range("J148").Value2 = "D"
$ws.Range("K148").Value2 = 1.285
$ws.Range("L148").Value2 = 4.5
$ws.Range("M148").Value2 = 9
$ws.Range("N148").Value2 = 1.5
$ws.Range("O148").Value2 = 3.8
$ws.Range("P148").Value2 = 5.5
$ws.Range("Q148").Value2 = -1
$ws.Range("R148").Value2 = 1.95
$ws.Range("S148").Value2 = 1.85
$ws.Range("T148").Value2 = 1.75
$ws.Range("U148").Value2 = 1.8
$ws.Range("V148").Value2 = 2
$ws.Range("W148").Value2 = -1
$ws.Range("X148").Value2 = 2.8
$ws.Range("Y148").Value2 = -1
$ws.Range("Z148").Value2 = -1
$ws.Range("AA148").Value2 = 0.8500000000000001
$ws.Range("AB148").Value2 = -1
$ws.Range("AC148").Value2 = 1

# Row 149 (previously row 147): renumber its running id (two new rows
# were inserted ahead of it), fill in the result + remaining odds
# columns that were missing, and refresh the PL columns with the
# computed figures.
$ws.Range("A149").Value2 = 147
$ws.Range("H149").Value2 = 2
$ws.Range("I149").Value2 = 1
$ws.Range("J149").Value2 = "H"
$ws.Range("W149").Value2 = 1.2
$ws.Range("X149").Value2 = -1
$ws.Range("Y149").Value2 = -1
$ws.Range("Z149").Value2 = 0.95
$ws.Range("AA149").Value2 = -1
$ws.Range("AB149").Value2 = 0.8500000000000001
$ws.Range("AC149").Value2 = -1
